$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.833.01'
$ws.Range('E2').Value = '  +3.23%  '

$ws.Range('D3').Value = '3.206.47'

$ws.Range('E4').Value = '  -0.14%  '

$ws.Range('D5').Value = '''603.04'
$ws.Range('E5').Value = '  +3.79%  '

$ws.Range('D6').Value = '''157.99'
$ws.Range('E6').Value = '  +6.42%  '

$ws.Range('D7').Value = '''1.00'
$ws.Range('E7').Value = '  -0.09%  '

$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').Value = '''0.555'
$ws.Range('E8').Value = '  +5.99%  '

$ws.Range('B9').Value = 'LidoStakedEther'
$ws.Range('C9').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D9').Value = '3.207.02'
$ws.Range('E9').Value = '  +2.15%  '

$ws.Range('D10').Value = '''0.161'
$ws.Range('E10').Value = '  +2.10%  '

$ws.Range('E11').Value = '  -1.48%  '

$ws.Range('E12').Value = '  +3.99%  '

$ws.Range('E13').Value = '  +1.76%  '

$ws.Range('D14').Value = '''39.38'
$ws.Range('E14').Value = '  +6.54%  '

$ws.Range('D15').Value = '3.735.90'
$ws.Range('E15').Value = '  +1.83%  '

$ws.Range('D16').Value = '66.850.08'
$ws.Range('E16').Value = '  +3.10%  '

$ws.Range('D17').Value = '''7.49'
$ws.Range('E17').Value = '  +5.03%  '

$ws.Range('D18').Value = '3.207.71'
$ws.Range('E18').Value = '  +1.88%  '

$ws.Range('E19').Value = '  +0.89%  '

$ws.Range('D20').Value = '''520.22'
$ws.Range('E20').Value = '  +4.00%  '

$ws.Range('D21').Value = '''15.47'
$ws.Range('E21').Value = '  +1.46%  '

$ws.Range('D22').Value = '''0.746'
$ws.Range('E22').Value = '  +4.90%  '

$ws.Range('D23').Value = '''8.24'
$ws.Range('E23').Value = '  +6.34%  '

$ws.Range('D24').Value = '''15.16'
$ws.Range('E24').Value = '  +1.02%  '

$ws.Range('D25').Value = '''85.47'
$ws.Range('E25').Value = '  +0.81%  '

$ws.Range('D26').Value = '''0.999'
$ws.Range('E26').Value = '  -0.14%  '

$ws.Range('D27').Value = '''9.44'
$ws.Range('E27').Value = '  +4.49%  '

$ws.Range('D28').Value = '''3.04'
$ws.Range('E28').Value = '  +4.39%  '

$ws.Range('D29').Value = '''2.44'
$ws.Range('E29').Value = '  +11.76%  '

$ws.Range('D30').Value = '''3.10'
$ws.Range('E30').Value = '  +11.45%  '

$ws.Range('D31').Value = '''6.98'
$ws.Range('E31').Value = '  +8.98%  '

$ws.Range('D32').Value = '''28.35'
$ws.Range('E32').Value = '  +3.00%  '

$ws.Range('E33').Value = '  +2.07%  '

$ws.Range('E34').Value = '  +0.23%  '

$ws.Range('D35').Value = '''6.61'
$ws.Range('E35').Value = '  +2.47%  '

$ws.Range('D36').Value = '''522.24'
$ws.Range('E36').Value = '  +10.91%  '

$ws.Range('D37').Value = '''55.01'
$ws.Range('E37').Value = '  +0.38%  '

$ws.Range('D38').Value = '''0.0910'
$ws.Range('E38').Value = '  +2.12%  '

$ws.Range('D39').Value = '''0.0427'
$ws.Range('E39').Value = '  +2.60%  '

$ws.Range('D40').Value = '''0.127'
$ws.Range('E40').Value = '  +9.24%  '

$ws.Range('B41').Value = 'dogwifhat'
$ws.Range('C41').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D41').Value = '''2.96'
$ws.Range('E41').Value = '  +1.85%  '

$ws.Range('B42').Value = 'TheGraph'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D42').Value = '''0.312'
$ws.Range('E42').Value = '  +11.26%  '

$ws.Range('B43').Value = 'Cosmos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D43').Value = '''8.96'
$ws.Range('E43').Value = '  +2.93%  '

$ws.Range('B44').Value = 'PEPE'
$ws.Range('C44').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D44').Value = '0.0₃0690'
$ws.Range('E44').Value = '  +14.35%  '

$ws.Range('D45').Value = '''2.53'
$ws.Range('E45').Value = '  +4.15%  '

$ws.Range('D46').Value = '2.903.18'
$ws.Range('E46').Value = '  -2.64%  '

$ws.Range('D47').Value = '''29.00'
$ws.Range('E47').Value = '  +2.86%  '

$ws.Range('D48').Value = '''2.43'
$ws.Range('E48').Value = '  +8.90%  '

$ws.Range('D49').Value = '''0.118'
$ws.Range('E49').Value = '  +3.86%  '

$ws.Range('E50').Value = '  +12.76%  '

$ws.Range('E51').Value = '  -0.01%  '
